$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.419.77"
$ws.Range("E2").Value = "  +2.63%  "
$ws.Range("D3").Value = "2.736.24"
$ws.Range("E3").Value = "  +2.58%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'115.70"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").Value = "'331.03"
$ws.Range("E6").Value = "  +1.23%  "
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.563"
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("D10").Value = "'41.54"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("D11").Value = "'20.35"
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("D12").Value = "'0.0828"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("E13").Value = "  +2.71%  "
$ws.Range("D14").Value = "'7.64"
$ws.Range("E14").Value = "  +3.26%  "
$ws.Range("D15").Value = "3.164.11"
$ws.Range("E15").Value = "  +2.69%  "
$ws.Range("D16").Value = "2.731.44"
$ws.Range("E16").Value = "  +1.65%  "
$ws.Range("D17").Value = "'0.883"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "51.333.67"
$ws.Range("E18").Value = "  +2.65%  "
$ws.Range("D19").Value = "'13.71"
$ws.Range("E19").Value = "  +2.72%  "
$ws.Range("D20").Value = "'3.04"
$ws.Range("E20").Value = "  +3.71%  "
$ws.Range("D21").Value = "'6.84"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").Value = "0.0₃0963"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "'286.54"
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("D24").Value = "'70.46"
$ws.Range("E24").Value = "  -3.27%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "'26.90"
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").Value = "'10.32"
$ws.Range("E28").Value = "  +2.61%  "
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("D30").Value = "'0.141"
$ws.Range("E30").Value = "  -1.43%  "
$ws.Range("D31").Value = "'35.74"
$ws.Range("E31").Value = "  -2.96%  "
$ws.Range("D32").Value = "'50.03"
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("E33").Value = "  +1.44%  "
$ws.Range("D34").Value = "'0.0827"
$ws.Range("E34").Value = "  +1.13%  "
$ws.Range("D35").Value = "'19.41"
$ws.Range("E35").Value = "  -2.05%  "
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").Value = "'5.04"
$ws.Range("E37").Value = "  -1.57%  "
$ws.Range("E38").Value = "  +1.23%  "
$ws.Range("D39").Value = "'3.23"
$ws.Range("E39").Value = "  +3.07%  "
$ws.Range("D40").Value = "'23.87"
$ws.Range("E40").Value = "  +7.54%  "
$ws.Range("D41").Value = "'129.34"
$ws.Range("E41").Value = "  +2.98%  "
$ws.Range("D42").Value = "'0.0350"
$ws.Range("E42").Value = "  +9.32%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "'0.113"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'2.29"
$ws.Range("E44").Value = "  +2.65%  "
$ws.Range("D45").Value = "'3.42"
$ws.Range("E45").Value = "  +2.08%  "
$ws.Range("D46").Value = "2.116.87"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").Value = "'2.24"
$ws.Range("E47").Value = "  +11.95%  "
$ws.Range("E48").Value = "  -2.38%  "
$ws.Range("D49").Value = "'5.51"
$ws.Range("E49").Value = "  +2.44%  "
$ws.Range("D50").Value = "'9.10"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("D51").Value = "'60.18"
$ws.Range("E51").Value = "  +1.00%  "

# Reset style to Normal for cells that needed a quote-prefix to stay text,
# so the style index matches the original (unstyled) cells.
foreach ($addr in @("D4","D5","D6","D8","D9","D10","D11","D12","D14","D17","D19","D20","D21","D23","D24","D26","D28","D30","D31","D32","D34","D35","D37","D39","D40","D41","D42","D43","D44","D45","D47","D49","D50","D51")) {
    $ws.Range($addr).Style = "Normal"
}
